# #439: migrate styles and classes, fixes from testing, add actions for accordion/tabs
#
# Rename two service names whose generated-GUID-style suffix changed
# (ca476 -> d97b7), and fix two rows where Status was swapped during
# testing (AppXSvc was actually Running, BDESVC was actually Stopped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 : AarSvc_ca476 -> AarSvc_d97b7
$ws.Cells.Item(2, 1).Value = "AarSvc_d97b7"

# Row 19 : BcastDVRUserService_ca476 -> BcastDVRUserService_d97b7
$ws.Cells.Item(19, 1).Value = "BcastDVRUserService_d97b7"

# Row 12 (AppXSvc) : Status Stopped -> Running
$ws.Cells.Item(12, 2).Value = "Running"

# Row 20 (BDESVC) : Status Running -> Stopped
$ws.Cells.Item(20, 2).Value = "Stopped"
